# OrangeHRMDDT/results/TestResults.xlsx -- "Edited code -- By Nirmala"
#
# Content changes:
#  1. On the "LoginValidData" sheet, the per-row "Pass" result column
#     (C2:C4) is removed entirely -- clear both the cached "Pass" text
#     and the green highlight style that went with it.
#  2. The workbook's active/selected tab moves from "LoginValidData" to
#     "LoginInvalidData".

$wb = $excel.ActiveWorkbook

$wsValid   = $wb.Worksheets.Item("LoginValidData")
$wsInvalid = $wb.Worksheets.Item("LoginInvalidData")

# Remove the stray "Pass" results column (values + styling) from the
# LoginValidData sheet.
$wsValid.Range("C2:C4").Clear()

# LoginInvalidData becomes the active sheet/tab.
$wsInvalid.Activate()
